# Fill in the manually-collected match scores for Week 7 (rows 46-49)
# and the advancing teams for the Prelim round (rows 53-54).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week 7 individual map scores (columns C/E, P/R, X/Z for rows 46-49) ---

$ws.Range("C46").Value = 6
$ws.Range("E46").Value = 4
$ws.Range("P46").Value = 13
$ws.Range("R46").Value = 6
$ws.Range("X46").Value = 7
$ws.Range("Z46").Value = 1

$ws.Range("C47").Value = 4
$ws.Range("E47").Value = 7
$ws.Range("P47").Value = 15
$ws.Range("R47").Value = 5
$ws.Range("X47").Value = 7
$ws.Range("Z47").Value = 7

$ws.Range("C48").Value = 5
$ws.Range("E48").Value = 9
$ws.Range("P48").Value = 6
$ws.Range("R48").Value = 6
$ws.Range("X48").Value = 8
$ws.Range("Z48").Value = 6

$ws.Range("C49").Value = 5
$ws.Range("E49").Value = 7
$ws.Range("P49").Value = 5
$ws.Range("R49").Value = 3
$ws.Range("X49").Value = 3
$ws.Range("Z49").Value = 7

# E48 loses its highlight shading once the score is filled in
$ws.Range("E48").Interior.Pattern = -4142

# --- Prelim round (row 53/54): teams that advanced out of Week 7 ---

$ws.Range("A53").Value = "Dexs Midnight Jukers"
$ws.Range("G53").Value = "Ballton Wanderers"
$ws.Range("N53").Value = "The Rickrollers"
$ws.Range("T53").Value = "Ballmeiras"

$ws.Range("A54").Value = "Ballton Wanderers"
$ws.Range("G54").Value = "Dexs Midnight Jukers"
$ws.Range("N54").Value = "Ballmeiras"
$ws.Range("T54").Value = "The Rickrollers"
